$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (this shifts existing rows 7..40 down to 8..41)
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly price record
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(7, 3).Value = "Los Lagos"
$ws.Cells.Item(7, 4).Value = 44532
$ws.Cells.Item(7, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 300000000
$ws.Cells.Item(7, 7).Value = "Espárragos"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 180
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = 1500
$ws.Cells.Item(7, 14).Value = "`$/kilo"
$ws.Cells.Item(7, 15).Value = "Provincia de Linares"
$ws.Cells.Item(7, 16).Value = 1500
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
